# "edit features of items"
# The "available" (F) column is converted from the old "yes"/"no" labels to
# the new, more descriptive "Available" / "Not Available" labels, and the
# per-row availability flags are refreshed at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$available = @{
    2  = "Not Available"
    3  = "Not Available"
    4  = "Not Available"
    5  = "Available"
    6  = "Available"
    7  = "Available"
    8  = "Available"
    9  = "Available"
    10 = "Not Available"
    11 = "Available"
    12 = "Available"
    13 = "Available"
    14 = "Available"
    15 = "Available"
    16 = "Available"
    17 = "Not Available"
    18 = "Available"
    19 = "Available"
    20 = "Available"
    21 = "Available"
    22 = "Available"
    23 = "Available"
    24 = "Available"
    25 = "Not Available"
}

foreach ($row in $available.Keys) {
    $ws.Range("F$row").Value = $available[$row]
}

# Let column F resize itself to fit the new, longer labels.
$ws.Columns("F").AutoFit()

# Restore the cursor position left behind by the edit.
$ws.Range("F27").Select() | Out-Null
